$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Recorded By" email lists in column G for the specified rows
# (values are unchanged as a set, only the comma-separated ordering changes)

$ws.Range("G3").Value2 = 'asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G4").Value2 = 'asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G5").Value2 = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range("G6").Value2 = 'asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Range("G7").Value2 = 'asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G8").Value2 = 'asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, backup@backdoor.com, eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G12").Value2 = 'System, salma.elgendy.std@med.asu.edu.eg'
$ws.Range("G13").Value2 = 'Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Range("G15").Value2 = 'wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G20").Value2 = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range("G24").Value2 = 'Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G25").Value2 = 'Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G26").Value2 = 'youstina.magdy@med.asu.edu.eg, gigi@med.asu.edu.eg'
$ws.Range("G27").Value2 = 'neveen.nashaat@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
$ws.Range("G30").Value2 = 'asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G31").Value2 = 'asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G32").Value2 = 'asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G33").Value2 = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range("G34").Value2 = 'asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Range("G35").Value2 = 'asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G36").Value2 = 'asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, backup@backdoor.com, eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G40").Value2 = 'System, salma.elgendy.std@med.asu.edu.eg'
$ws.Range("G41").Value2 = 'Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Range("G43").Value2 = 'wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G48").Value2 = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range("G52").Value2 = 'Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G53").Value2 = 'Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G54").Value2 = 'youstina.magdy@med.asu.edu.eg, gigi@med.asu.edu.eg'
$ws.Range("G55").Value2 = 'neveen.nashaat@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
